$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1683.8
$ws.Range("J46").Value = 2139.6667
$ws.Range("L46").Value = 6419.000100000001
$ws.Range("N46").Value = -6657.000100000001
$ws.Range("H60").Value = 1683.8
$ws.Range("J60").Value = 2139.6667
$ws.Range("L60").Value = 6419.000100000001
$ws.Range("N60").Value = -7387.000100000001
$ws.Range("H127").Value = 229722.52
$ws.Range("I127").Value = 352.85715
$ws.Range("J127").Value = 269862.2
$ws.Range("K127").Value = 1058.57145
$ws.Range("L127").Value = 809586.6000000001
$ws.Range("M127").Value = 3901.42855
$ws.Range("N127").Value = -819506.6000000001
$ws.Range("H129").Value = 1996.9811
$ws.Range("J129").Value = 2321.5908
$ws.Range("L129").Value = 6964.7724
$ws.Range("N129").Value = -16964.7724
$ws.Range("H132").Value = 3849337.5
$ws.Range("I132").Value = 3008.1035
$ws.Range("J132").Value = 35718924
$ws.Range("K132").Value = 9024.3105
$ws.Range("L132").Value = 107156772
$ws.Range("M132").Value = -6494.3105
$ws.Range("N132").Value = -107161832
$ws.Range("H137").Value = 8000910.5
$ws.Range("I137").Value = 656
$ws.Range("J137").Value = 20001292
$ws.Range("K137").Value = 1968
$ws.Range("L137").Value = 60003876
$ws.Range("M137").Value = 582
$ws.Range("N137").Value = -60008976
$ws.Range("H138").Value = 9525575
$ws.Range("I138").Value = 13890124
$ws.Range("J138").Value = 2922.7273
$ws.Range("K138").Value = 41670372
$ws.Range("L138").Value = 8768.1819
$ws.Range("M138").Value = -41665232
$ws.Range("N138").Value = -19048.1819
$ws.Range("H141").Value = 1682.037
$ws.Range("I141").Value = 1564.5454
$ws.Range("J141").Value = 2199
$ws.Range("K141").Value = 4693.6362
$ws.Range("L141").Value = 6597
$ws.Range("M141").Value = 486.3638000000001
$ws.Range("N141").Value = -16957

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10374.694
$ws.Range("I32").Value = 10979.17
$ws.Range("K32").Value = 10979.17
$ws.Range("M32").Value = -10692.17
$ws.Range("H61").Value = 10871089
$ws.Range("I61").Value = 11906297
$ws.Range("K61").Value = 11906297
$ws.Range("M61").Value = -11906085
$ws.Range("H74").Value = 10640693
$ws.Range("I74").Value = 13890626
$ws.Range("J74").Value = 4546.727
$ws.Range("K74").Value = 13890626
$ws.Range("L74").Value = 4546.727
$ws.Range("M74").Value = -13889752
$ws.Range("N74").Value = -6294.727
$ws.Range("H77").Value = 10640693
$ws.Range("I77").Value = 13890626
$ws.Range("J77").Value = 4546.727
$ws.Range("K77").Value = 69453130
$ws.Range("L77").Value = 22733.635
$ws.Range("M77").Value = -69448762
$ws.Range("N77").Value = -31469.635
$ws.Range("H122").Value = 14518.556
$ws.Range("I122").Value = 23342.4
$ws.Range("J122").Value = 3488.75
$ws.Range("K122").Value = 70027.20000000001
$ws.Range("L122").Value = 10466.25
$ws.Range("M122").Value = -67577.20000000001
$ws.Range("N122").Value = -15366.25
$ws.Range("H136").Value = 10871089
$ws.Range("I136").Value = 11906297
$ws.Range("K136").Value = 35718891
$ws.Range("M136").Value = -35716341

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1565
$ws.Range("I99").Value = 1599.6
$ws.Range("J99").Value = 1478.5
$ws.Range("K99").Value = 1599.6
$ws.Range("L99").Value = 1478.5
$ws.Range("M99").Value = -101.5999999999999
$ws.Range("N99").Value = -4474.5
$ws.Range("H108").Value = 12592
$ws.Range("I108").Value = 4500
$ws.Range("J108").Value = 20684
$ws.Range("K108").Value = 4500
$ws.Range("L108").Value = 20684
$ws.Range("M108").Value = -660
$ws.Range("N108").Value = -28364
$ws.Range("H126").Value = 1565
$ws.Range("I126").Value = 1599.6
$ws.Range("J126").Value = 1478.5
$ws.Range("K126").Value = 4798.799999999999
$ws.Range("L126").Value = 4435.5
$ws.Range("M126").Value = -2328.799999999999
$ws.Range("N126").Value = -9375.5
$ws.Range("H132").Value = 10418364
$ws.Range("I132").Value = 12501568
$ws.Range("J132").Value = 2341
$ws.Range("K132").Value = 37504704
$ws.Range("L132").Value = 7023
$ws.Range("M132").Value = -37502174
$ws.Range("N132").Value = -12083
$ws.Range("H134").Value = 351308.97
$ws.Range("I134").Value = 1150.8889
$ws.Range("K134").Value = 3452.6667
$ws.Range("M134").Value = -917.6666999999998
$ws.Range("H140").Value = 24395.8
$ws.Range("J140").Value = 24395.8
$ws.Range("L140").Value = 24395.8
$ws.Range("N140").Value = -34755.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5100
$ws.Range("I80").Value = 10000
$ws.Range("J80").Value = 2650
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 7950
$ws.Range("M80").Value = -29064
$ws.Range("N80").Value = -9822
$ws.Range("H83").Value = 5100
$ws.Range("I83").Value = 10000
$ws.Range("J83").Value = 2650
$ws.Range("K83").Value = 90000
$ws.Range("L83").Value = 23850
$ws.Range("M83").Value = -85320
$ws.Range("N83").Value = -33210
$ws.Range("H109").Value = 3586.3635
$ws.Range("I109").Value = 916.6667
$ws.Range("J109").Value = 4587.5
$ws.Range("K109").Value = 2750.0001
$ws.Range("L109").Value = 13762.5
$ws.Range("M109").Value = -1710.0001
$ws.Range("N109").Value = -15842.5
$ws.Range("H136").Value = 2341
$ws.Range("I136").Value = 1562.5
$ws.Range("J136").Value = 3586.6
$ws.Range("K136").Value = 4687.5
$ws.Range("L136").Value = 10759.8
$ws.Range("M136").Value = 412.5
$ws.Range("N136").Value = -20959.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4380.4736
$ws.Range("I126").Value = 2387.5
$ws.Range("J126").Value = 5829.909
$ws.Range("K126").Value = 7162.5
$ws.Range("L126").Value = 17489.727
$ws.Range("M126").Value = -4692.5
$ws.Range("N126").Value = -22429.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6309.758
$ws.Range("I7").Value = 7373.7144
$ws.Range("K7").Value = 7373.7144
$ws.Range("M7").Value = -7261.7144
$ws.Range("H9").Value = 455.7143
$ws.Range("I9").Value = 365
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 365
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = -141
$ws.Range("N9").Value = -1448
$ws.Range("H126").Value = 6309.758
$ws.Range("I126").Value = 7373.7144
$ws.Range("K126").Value = 22121.1432
$ws.Range("M126").Value = -19651.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 786.2545
$ws.Range("I136").Value = 619.3542
$ws.Range("J136").Value = 1930.7142
$ws.Range("K136").Value = 1858.0626
$ws.Range("L136").Value = 5792.142599999999
$ws.Range("M136").Value = 691.9374
$ws.Range("N136").Value = -10892.1426
